$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.1281536994595948
$ws.Range("C2").Value = 0.2797649995624122
$ws.Range("D2").Value = 0.3038831578311377
$ws.Range("E2").Value = 0.5512559821273033
$ws.Range("F2").Value = 0.5413836472969831

# Row 3
$ws.Range("B3").Value = 0.3997486404630072
$ws.Range("C3").Value = 0.6318827147352428
$ws.Range("D3").Value = 1.377332145493524
$ws.Range("E3").Value = 1.1735979488281
$ws.Range("F3").Value = 1.114398417685673

# Row 4
$ws.Range("B4").Value = 0.6166457405447743
$ws.Range("C4").Value = 1.030434176077713
$ws.Range("D4").Value = 4.226220055432564
$ws.Range("E4").Value = 2.055777238767023
$ws.Range("F4").Value = 1.981024288269795

# Row 5
$ws.Range("B5").Value = 0.3564886474955841
$ws.Range("C5").Value = 1.153557125445555
$ws.Range("D5").Value = 5.502783603765931
$ws.Range("E5").Value = 2.345801271157881
$ws.Range("F5").Value = 2.34258258903974
$ws.Range("G5").Value = 49

# Row 6
$ws.Range("B6").Value = 0.3339033228567314
$ws.Range("C6").Value = 1.171432943203542
$ws.Range("D6").Value = 5.596202172583759
$ws.Range("E6").Value = 2.365629339643842
$ws.Range("F6").Value = 2.366729117756491
$ws.Range("G6").Value = 48

# Row 7
$ws.Range("B7").Value = 0.3166843790088169
$ws.Range("C7").Value = 1.268286245602126
$ws.Range("D7").Value = 6.364584808664379
$ws.Range("E7").Value = 2.522812876268151
$ws.Range("F7").Value = 2.535576010015387
$ws.Range("G7").Value = 39

# Row 8
$ws.Range("B8").Value = 0.2837060048817768
$ws.Range("C8").Value = 1.270563093429218
$ws.Range("D8").Value = 6.590865384712885
$ws.Range("E8").Value = 2.567268078076944
$ws.Range("F8").Value = 2.585794346692225
$ws.Range("G8").Value = 38

# Row 9
$ws.Range("B9").Value = 0.06262643962625501
$ws.Range("C9").Value = 1.720805088009817
$ws.Range("D9").Value = 10.81057706780112
$ws.Range("E9").Value = 3.287944200834485
$ws.Range("F9").Value = 3.368529018236877
$ws.Range("G9").Value = 21

# Row 10
$ws.Range("B10").Value = -0.5304237684279632
$ws.Range("C10").Value = 1.353147886947694
$ws.Range("D10").Value = 6.53100329151635
$ws.Range("E10").Value = 2.555582769451295
$ws.Range("F10").Value = 2.594300777942687
$ws.Range("G10").Value = 14

# Row 11
$ws.Range("B11").Value = 0.7489990067474906
$ws.Range("C11").Value = 0.7489990067474906
$ws.Range("D11").Value = 0.8040835518222755
$ws.Range("E11").Value = 0.8967070602054361
$ws.Range("F11").Value = 0.5512304868582063
